$wb = $excel.ActiveWorkbook

# --- Add the new "AddVoucher" worksheet as the last tab (after TrackingOrder) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AddVoucher"

# --- Header row ---
$ws.Range("A1").Value = "Voucher"
$ws.Range("B1").Value = "Output"

# --- Voucher test cases ---
$ws.Range("A2").Value = "LAZADA123"
$ws.Range("B2").Value = "Sorry, this voucher is not valid. Please check for any typing errors."

$ws.Range("A3").Value = "GIATOT123"
$ws.Range("B3").Value = "Sorry, this voucher is not valid. Please check for any typing errors."

$ws.Range("A4").Value = "LAZADANEWYEAR123"
$ws.Range("B4").Value = "Sorry, this voucher is not valid. Please check for any typing errors."

# --- Widen column B so the long message is readable ---
$ws.Columns.Item(2).ColumnWidth = 48.75

# --- Leave the selection where the author left it when the sheet was saved ---
$ws.Range("B8").Select()
